$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(0.6408044419705359, -265.1326953808737, -2.699866219999995, 0.7914892148354087, 0.5109184350519644, 0.2132339996864685, 157.9878643119225, 0.3186216174995813, 0.02101025240066041, 0.1698159349501208, 0.2405142646481177, 0.4617726710043249, 0.2163006006629874, 0.4814312896101858, 29.09073025240775, 44.93611597569436)

for ($row = 2; $row -le 26; $row++) {
    for ($col = 2; $col -le 17; $col++) {
        $ws.Cells.Item($row, $col).Value = $values[$col - 2]
    }
}
